$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price/Volume columns in this sheet are text-formatted strings (e.g.
# "69.525.88" thousand-separated prices, or values with a leading/trailing
# space and a "%" sign). Force column D cells that would otherwise be
# auto-detected as numbers to stay as text, matching the source data.

$ws.Range("D2").Value = "69.525.88"
$ws.Range("E2").Value = "  +1.73%  "
$ws.Range("D3").Value = "3.909.54"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.18"
$ws.Range("E5").Value = "  +9.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.62"
$ws.Range("E6").Value = "  -1.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.614"
$ws.Range("E7").Value = "  -1.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.727"
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.171"
$ws.Range("E10").Value = "  +2.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000330"
$ws.Range("E11").Value = "  -4.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.31"
$ws.Range("E12").Value = "  -2.00%  "
$ws.Range("D13").Value = "4.545.14"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.24"
$ws.Range("E14").Value = "  -4.99%  "
$ws.Range("D15").Value = "3.913.39"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("E16").Value = "  +8.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.135"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.85"
$ws.Range("E18").Value = "  -3.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.85"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").Value = "69.452.69"
$ws.Range("E20").Value = "  +1.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "433.79"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.34"
$ws.Range("E22").Value = "  -4.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.23"
$ws.Range("E23").Value = "  -6.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.72"
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("E25").Value = "  +11.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.56"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.64"
$ws.Range("E27").Value = "  -5.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.09"
$ws.Range("E28").Value = "  -5.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "692.40"
$ws.Range("E29").Value = "  -3.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.15"
$ws.Range("E30").Value = "  -4.60%  "
$ws.Range("E31").Value = "  -3.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.83"
$ws.Range("E32").Value = "  -3.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "68.28"
$ws.Range("E33").Value = "  +12.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.443"
$ws.Range("E34").Value = "  +12.08%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "40.16"
$ws.Range("E35").Value = "  -3.78%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.87"
$ws.Range("E36").Value = "  -4.84%  "
$ws.Range("D37").Value = "0.0₃0831"
$ws.Range("E37").Value = "  -6.07%  "
$ws.Range("E38").Value = "  +3.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0479"
$ws.Range("E41").Value = "  -3.81%  "
$ws.Range("E42").Value = "  +3.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.75"
$ws.Range("E43").Value = "  -8.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.94"
$ws.Range("E44").Value = "  -5.77%  "
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.05"
$ws.Range("E47").Value = "  +7.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.29"
$ws.Range("E48").Value = "  -3.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "143.50"
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("D50").Value = "0.0₆0338"
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.05"
$ws.Range("E51").Value = "  -4.53%  "
